$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.039.83'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.566.13'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.88'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.90'
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.59'
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.958.92'
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.606.61'
$ws.Range('E15').Value = '  +3.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.05'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.844'
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.090.31'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.84'
$ws.Range('E19').Value = '  +4.11%  '
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.53'
$ws.Range('E22').Value = '  -0.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.01'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.96'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +3.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.73'
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.17'
$ws.Range('E29').Value = '  +2.06%  '
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('E31').Value = '  -4.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.79'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.40'
$ws.Range('E33').Value = '  +3.82%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('E35').Value = '  +2.52%  '
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.06'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('E39').Value = '  +5.02%  '
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.42'
$ws.Range('E41').Value = '  -5.32%  '
$ws.Range('E42').Value = '  +2.83%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.28'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.000.57'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('E47').Value = '  +2.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '83.70'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.812.13'
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.28'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('E51').Value = '  +2.75%  '
